{"js": "// Each entry is [oldText, newText] for one table-cell answer that changed\n// between the \"before\" and \"after\" revisions of the worksheet.\nconst replacements = [\n  [\"254\u00f73=84, 2\", \"447\u00f78=55, 7\"],\n  [\"778\u00f77=111, 1\", \"555\u00f79=61, 6\"],\n  [\"809\u00f73=269, 2\", \"798\u00f77=114, 0\"],\n  [\"566\u00f72=283, 0\", \"583\u00f73=194, 1\"],\n  [\"828\u00f75=165, 3\", \"375\u00f73=125, 0\"],\n  [\"564\u00f79=62, 6\", \"222\u00f79=24, 6\"],\n  [\"408\u00f72=204, 0\", \"280\u00f73=93, 1\"],\n  [\"138\u00f76=23, 0\", \"944\u00f72=472, 0\"],\n  [\"694\u00f79=77, 1\", \"449\u00f74=112, 1\"],\n  [\"332\u00f78=41, 4\", \"732\u00f77=104, 4\"],\n  [\"561\u00f72=280, 1\", \"639\u00f76=106, 3\"],\n  [\"662\u00f79=73, 5\", \"588\u00f72=294, 0\"],\n  [\"114\u00f77=16, 2\", \"739\u00f74=184, 3\"],\n  [\"733\u00f76=122, 1\", \"741\u00f74=185, 1\"],\n  [\"307\u00f74=76, 3\", \"925\u00f76=154, 1\"],\n  [\"882\u00f72=441, 0\", \"724\u00f77=103, 3\"],\n  [\"401\u00f78=50, 1\", \"853\u00f77=121, 6\"],\n  [\"215\u00f76=35, 5\", \"574\u00f72=287, 0\"],\n  [\"117\u00f79=13, 0\", \"194\u00f78=24, 2\"],\n  [\"147\u00f73=49, 0\", \"512\u00f72=256, 0\"],\n  [\"281\u00f74=70, 1\", \"104\u00f76=17, 2\"],\n  [\"955\u00f73=318, 1\", \"891\u00f74=222, 3\"],\n  [\"902\u00f72=451, 0\", \"487\u00f78=60, 7\"],\n  [\"587\u00f77=83, 6\", \"500\u00f79=55, 5\"],\n  [\"726\u00f72=363, 0\", \"976\u00f78=122, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  // Every answer string is unique in the document, so a plain, case-sensitive\n  // search reliably finds the single table cell that needs updating.\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Each pair is (oldText, newText) for one table-cell answer that changed\n# between the \"before\" and \"after\" revisions of the worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"254\u00f73=84, 2\", \"447\u00f78=55, 7\"),\n    @(\"778\u00f77=111, 1\", \"555\u00f79=61, 6\"),\n    @(\"809\u00f73=269, 2\", \"798\u00f77=114, 0\"),\n    @(\"566\u00f72=283, 0\", \"583\u00f73=194, 1\"),\n    @(\"828\u00f75=165, 3\", \"375\u00f73=125, 0\"),\n    @(\"564\u00f79=62, 6\", \"222\u00f79=24, 6\"),\n    @(\"408\u00f72=204, 0\", \"280\u00f73=93, 1\"),\n    @(\"138\u00f76=23, 0\", \"944\u00f72=472, 0\"),\n    @(\"694\u00f79=77, 1\", \"449\u00f74=112, 1\"),\n    @(\"332\u00f78=41, 4\", \"732\u00f77=104, 4\"),\n    @(\"561\u00f72=280, 1\", \"639\u00f76=106, 3\"),\n    @(\"662\u00f79=73, 5\", \"588\u00f72=294, 0\"),\n    @(\"114\u00f77=16, 2\", \"739\u00f74=184, 3\"),\n    @(\"733\u00f76=122, 1\", \"741\u00f74=185, 1\"),\n    @(\"307\u00f74=76, 3\", \"925\u00f76=154, 1\"),\n    @(\"882\u00f72=441, 0\", \"724\u00f77=103, 3\"),\n    @(\"401\u00f78=50, 1\", \"853\u00f77=121, 6\"),\n    @(\"215\u00f76=35, 5\", \"574\u00f72=287, 0\"),\n    @(\"117\u00f79=13, 0\", \"194\u00f78=24, 2\"),\n    @(\"147\u00f73=49, 0\", \"512\u00f72=256, 0\"),\n    @(\"281\u00f74=70, 1\", \"104\u00f76=17, 2\"),\n    @(\"955\u00f73=318, 1\", \"891\u00f74=222, 3\"),\n    @(\"902\u00f72=451, 0\", \"487\u00f78=60, 7\"),\n    @(\"587\u00f77=83, 6\", \"500\u00f79=55, 5\"),\n    @(\"726\u00f72=363, 0\", \"976\u00f78=122, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # Every answer string is unique in the document, so a plain,\n    # case-sensitive Find/Replace reliably hits the single table cell\n    # that needs updating.\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
